$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.577.68'
$ws.Range("E2").Value = '  +5.93%  '
$ws.Range("D3").Value = '2.739.29'
$ws.Range("E3").Value = '  +4.80%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'116.11"
$ws.Range("E5").Value = '  +5.70%  '
$ws.Range("D6").Value = "'333.27"
$ws.Range("E6").Value = '  +3.57%  '
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = "'0.568"
$ws.Range("E9").Value = '  +5.27%  '
$ws.Range("D10").Value = "'41.30"
$ws.Range("E10").Value = '  +4.95%  '
$ws.Range("D11").Value = "'0.0852"
$ws.Range("E11").Value = '  +5.48%  '
$ws.Range("D12").Value = "'20.08"
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("E13").Value = '  +2.81%  '
$ws.Range("D14").Value = "'7.53"
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("D15").Value = '3.179.36'
$ws.Range("E15").Value = '  +5.30%  '
$ws.Range("D16").Value = '2.751.26'
$ws.Range("E16").Value = '  +4.17%  '
$ws.Range("D17").Value = "'0.874"
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").Value = '51.530.73'
$ws.Range("E18").Value = '  +5.99%  '
$ws.Range("D19").Value = "'3.11"
$ws.Range("E19").Value = '  +5.74%  '
$ws.Range("D20").Value = "'13.41"
$ws.Range("E20").Value = '  +4.57%  '
$ws.Range("D21").Value = "'6.81"
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").Value = '0.0₃0971'
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").Value = "'278.70"
$ws.Range("E23").Value = '  +3.51%  '
$ws.Range("D24").Value = "'69.27"
$ws.Range("E24").Value = '  +1.02%  '
$ws.Range("D25").Value = "'2.64"
$ws.Range("E25").Value = '  +4.48%  '
$ws.Range("D26").Value = "'26.64"
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "'10.15"
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").Value = "'0.139"
$ws.Range("E30").Value = '  +2.66%  '
$ws.Range("D31").Value = "'34.82"
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").Value = "'49.95"
$ws.Range("E32").Value = '  +1.42%  '
$ws.Range("D33").Value = "'5.53"
$ws.Range("E33").Value = '  +1.56%  '
$ws.Range("D34").Value = "'0.0816"
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = "'18.90"
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").Value = "'4.95"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("E39").Value = '  +1.40%  '
$ws.Range("D40").Value = "'127.57"
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("D41").Value = "'0.0344"
$ws.Range("E41").Value = '  +8.86%  '
$ws.Range("E42").Value = '  +2.17%  '
$ws.Range("E43").Value = '  +2.28%  '
$ws.Range("D44").Value = "'2.27"
$ws.Range("E44").Value = '  +6.85%  '
$ws.Range("E45").Value = '  +12.23%  '
$ws.Range("D46").Value = '2.085.67'
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = '  +2.55%  '
$ws.Range("E48").Value = '  +3.05%  '
$ws.Range("D49").Value = "'5.51"
$ws.Range("E49").Value = '  +7.18%  '
$ws.Range("D50").Value = "'8.91"
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("D51").Value = "'59.57"
$ws.Range("E51").Value = '  +1.91%  '
